# Added Test Data for UK Market
# Duplicate the "Netherlands" sheet (last tab) to create a new "UK" tab,
# placed after it, then update the new sheet's "User Story" input cell
# (B4) with the UK-specific test-case reference.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("Netherlands")

# Copy places the new sheet immediately after the given sheet and makes
# it the active sheet - matching Excel's "Move or Copy... (Create a copy)".
$sourceSheet.Copy([System.Reflection.Missing]::Value, $sourceSheet)

$ukSheet = $wb.Worksheets.Item($sourceSheet.Index + 1)
$ukSheet.Name = "UK"

$ukSheet.Range("B4").Value = "NGC-2741/T3399"
